# Generate Report for Handback
# Adds a new handback row (a4490a95-e742-4852-8490-6174db646455) to the
# Overview sheet and to each language sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": append row 3
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A3").Value = "a4490a95-e742-4852-8490-6174db646455.md"
$ws1.Range("B3").Value = "e2e\a4490a95-e742-4852-8490-6174db646455.md"
$ws1.Range("C3").Value = ".md"
$ws1.Range("E3").Value = "Handed back: in sync with en-US"
$ws1.Range("F3").Value = "Handed back: in sync with en-US"
$ws1.Range("G3").Value = "2016-10-10 09:27:06"
$ws1.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws1.Range("B3").Font.Underline = $true
$ws1.Range("B3").Font.Color = 15570276

$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48cde75648e9a4cdadf07ed818a7ab23a861da49/e2e/a4490a95-e742-4852-8490-6174db646455.md", "", "", "e2e\a4490a95-e742-4852-8490-6174db646455.md")

$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn": append row 3
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A3").Value = "a4490a95-e742-4852-8490-6174db646455.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = "True"
$ws2.Range("G3").Value = "a4490a95-e742-4852-8490-6174db646455.a67382c08b6a095cb991297b872e1c4f9e9c2baa.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-10-10 09:26:56"
$ws2.Range("I3").Value = "a4490a95-e742-4852-8490-6174db646455.md"
$ws2.Range("J3").Value = "a4490a95-e742-4852-8490-6174db646455.a67382c08b6a095cb991297b872e1c4f9e9c2baa.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-10-10 09:27:42"
$ws2.Range("L3").Value = ""
$ws2.Range("M3").Value = "True"
$ws2.Range("N3").Value = ""
$ws2.Range("O3").Value = "False"
$ws2.Range("P3").Value = ""

$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Range("A3").Font.Underline = $true
$ws2.Range("A3").Font.Color = 15570276
$ws2.Range("I3").Font.Underline = $true
$ws2.Range("I3").Font.Color = 15570276

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/656bac3ae93bd3e349458c2c626306fd54f12112/e2e/a4490a95-e742-4852-8490-6174db646455.md", "", "", "a4490a95-e742-4852-8490-6174db646455.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/656bac3ae93bd3e349458c2c626306fd54f12112/e2e/a4490a95-e742-4852-8490-6174db646455.md", "", "", "a4490a95-e742-4852-8490-6174db646455.md")

$lo2 = $ws2.ListObjects.Item(1)
$lo2.Resize($ws2.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet "de-de": append row 3
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A3").Value = "a4490a95-e742-4852-8490-6174db646455.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = "True"
$ws3.Range("G3").Value = "a4490a95-e742-4852-8490-6174db646455.a67382c08b6a095cb991297b872e1c4f9e9c2baa.de-de.xlf"
$ws3.Range("H3").Value = "2016-10-10 09:27:06"
$ws3.Range("I3").Value = "a4490a95-e742-4852-8490-6174db646455.md"
$ws3.Range("J3").Value = "a4490a95-e742-4852-8490-6174db646455.a67382c08b6a095cb991297b872e1c4f9e9c2baa.de-de.xlf"
$ws3.Range("K3").Value = "2016-10-10 09:27:57"
$ws3.Range("L3").Value = ""
$ws3.Range("M3").Value = "True"
$ws3.Range("N3").Value = ""
$ws3.Range("O3").Value = "False"
$ws3.Range("P3").Value = ""

$ws3.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws3.Range("A3").Font.Underline = $true
$ws3.Range("A3").Font.Color = 15570276
$ws3.Range("I3").Font.Underline = $true
$ws3.Range("I3").Font.Color = 15570276

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/82ffc4b48fe421eefb2fc2a0355db11041885ffd/e2e/a4490a95-e742-4852-8490-6174db646455.md", "", "", "a4490a95-e742-4852-8490-6174db646455.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/82ffc4b48fe421eefb2fc2a0355db11041885ffd/e2e/a4490a95-e742-4852-8490-6174db646455.md", "", "", "a4490a95-e742-4852-8490-6174db646455.md")

$lo3 = $ws3.ListObjects.Item(1)
$lo3.Resize($ws3.Range("A1:P3"))

Write-Output "Generate Report for Handback: done"
